# Update "All Orders" sheet: order #9 (row 11) status changes from NEW to
# CANCELLED, and a cancel reason of "test" is recorded.
$wbAll = $excel.ActiveWorkbook
$wsOrders = $wbAll.Worksheets.Item("All Orders")
$wsOrders.Range("H11").Value = "CANCELLED"
$wsOrders.Range("M11").Value = "test"

# Update "Daily Summary" sheet for 2026-01-13 (row 4): one more cancellation,
# and revenue/pending drop to 0 since that order's amount no longer counts.
$wsSummary = $wbAll.Worksheets.Item("Daily Summary")
$wsSummary.Range("D4").Value = 13
$wsSummary.Range("E4").Value = 0
$wsSummary.Range("G4").Value = 0
